# Title page fix ("titlesec bug in Ubuntu 16.04"):
# Prepend a bold "Master " run in front of the existing "Thesis Report"
# run so the title reads "Master Thesis Report", matching the formatting
# (Cambria / majorHAnsi theme font, bold, en-US) already used by the
# neighbouring run.

$d = $word.ActiveDocument

# Locate the start of the (unique) "Thesis Report" run.
$found = $d.Content
$found.Find.Execute("Thesis Report", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$insertionPoint = $found.Start

# Insert the new text as its own run immediately before "Thesis Report".
$ins = $d.Range($insertionPoint, $insertionPoint)
$ins.InsertBefore("Master ")

# The newly typed text inherits identical run formatting to "Thesis Report",
# so the engine would otherwise coalesce the two into a single <w:r>. Drop a
# throwaway bookmark at the boundary to force the run split, then remove the
# bookmark again so no trace of it is left in the saved document.
$boundary = $insertionPoint + 7
$splitMark = $d.Range($boundary, $boundary)
$d.Bookmarks.Add("zzz_run_split_marker", $splitMark)
$d.Bookmarks("zzz_run_split_marker").Delete()
